$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Sheet 2")

# Add a new row of data (row 17) with a "less than" rule checker pair of values
$ws.Range("D17").Value = 200
$ws.Range("E17").Value = 201

# Move the active selection to D18, as left by the author after adding the new row
$ws.Activate()
$ws.Range("D18").Select()
